$wb = $excel.ActiveWorkbook

# --- Add the new "Reg" worksheet as the last sheet in the workbook ---
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$ws.Name = "Reg"

# Note: values are entered in this particular order (A2 "Second" before B1
# "Val") so that the new shared-string table entries come out in the same
# order as the target workbook (Second=20, Val=21).
$ws.Range("A2").Value = "Second"

$ws.Range("B1").Value = "Val"
$ws.Range("B1").Font.Bold = $true
$ws.Range("B1").HorizontalAlignment = -4152   # xlRight

$ws.Range("B2").Value = 15
$ws.Range("B2").HorizontalAlignment = -4152   # xlRight

$ws.Range("C2").Formula = "=_xlfn.BITAND(B2,15)"
$ws.Range("C2").HorizontalAlignment = -4152   # xlRight

$ws.Range("B3").Value = 166
$ws.Range("B3").HorizontalAlignment = -4152   # xlRight

$ws.Range("C3").Formula = "=DEC2BIN(B3, 8)"
$ws.Range("C3").HorizontalAlignment = -4152   # xlRight

$ws.PageSetup.Orientation = 1   # xlPortrait

# --- Activate the new sheet, freeze the header row and set the selection ---
$ws.Activate() | Out-Null
$ws.Range("A2").Select() | Out-Null
$excel.ActiveWindow.FreezePanes = $true
$ws.Range("C4").Select() | Out-Null
